$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (shifts existing B.. columns to D..)
$ws.Range("B:C").Insert()

# Row 1 - headers (new dates)
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"

# Row 2
$ws.Range("B2").Value = "Hold       ($10.21)"
$ws.Range("C2").Value = "Hold       ($10.09)"

# Row 3
$ws.Range("B3").Value = "Hold       ($25.08)"
$ws.Range("C3").Value = "Hold       ($25.13)"

# Row 4
$ws.Range("B4").Value = "Hold       ($17.14)"
$ws.Range("C4").Value = "Hold       ($17.45)"

# Row 5
$ws.Range("B5").Value = "UN         (0)"
$ws.Range("C5").Value = "UN         (0)"

# Row 6
$ws.Range("B6").Value = "Buy        ($64.41)"
$ws.Range("C6").Value = "Buy        ($64.03)"
